$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "Alessandro GALVAGNI"
$ws.Range("B28").Value = "Thomas Debiasi | MAI UNA GIOIA"
$ws.Range("C28").Value = "Michele Merighi | Clitoriders"
$ws.Range("D28").Value = "Davide Raffaelli | MediaserT"
$ws.Range("E28").Value = "MATTEO BRIGO | Pinguini Trentini"
$ws.Range("F28").Value = "Daniele Ruzzenenti | Demobusters"
